$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell, preventing Excel's
# automatic number/date inference (so values like "1234" or "01/01/2000"
# stay as text, matching the source application's export behaviour), and
# then reset the cell style back to Normal so no extra numFmt/quote-prefix
# style lingers on the cell itself.
function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 3 values (columns A-W)
$row3 = @(
    "01/01/2000",
    "wefwfwf",
    "AgroIndústria",
    "Codigestão",
    "Abatedouro de aves",
    "Água residual",
    "wefwf",
    "wefwfwfw",
    "", "", "", "", "", "", "", "", "", "", "", "", "", "", ""
)

# Row 4 values (columns A-W)
$row4 = @(
    "01/01/2000",
    "ragagegtrbbr",
    "Resíduo urbano",
    "Resíduos alimentares",
    "Cebola",
    "Amostra suco",
    "wergwergwerg",
    "wergrtrtrt",
    "1",
    "",
    "",
    "2",
    "",
    "",
    "3",
    "",
    "",
    "1234",
    "",
    "",
    "1234124",
    "",
    ""
)

for ($i = 0; $i -lt $row3.Length; $i++) {
    Set-TextValue $ws.Cells.Item(3, $i + 1) $row3[$i]
}

for ($i = 0; $i -lt $row4.Length; $i++) {
    Set-TextValue $ws.Cells.Item(4, $i + 1) $row4[$i]
}
